$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new empty row above row 1, shifting the header row and the
# existing transaction rows down by one (rows 1-4 -> rows 2-5).
$ws.Rows.Item(1).Insert()

# Add the new row of data (1, 2, 3) right after the existing table.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 3
